$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.645.24"
$ws.Range("D3").Value = "3.547.18"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'580.20"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").Value = "'186.76"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("D8").Value = "3.535.71"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D10").Value = "'0.218"
$ws.Range("E10").Value = "  +18.21%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "'54.21"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "'0.0000316"
$ws.Range("E13").Value = "  +5.21%  "
$ws.Range("D14").Value = "'9.46"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "4.114.15"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "70.690.19"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "'19.08"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "'12.72"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'578.54"
$ws.Range("E19").Value = "  +6.72%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.504.30"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "'17.65"
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("D24").Value = "'4.55"
$ws.Range("E24").Value = "  +2.76%  "
$ws.Range("D25").Value = "'4.89"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("D26").Value = "'94.31"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "'11.15"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").Value = "'9.29"
$ws.Range("D30").Value = "'32.60"
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("D31").Value = "'7.17"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").Value = "'12.27"
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").Value = "'3.75"
$ws.Range("E35").Value = "  +22.04%  "
$ws.Range("E36").Value = "  +6.45%  "
$ws.Range("D37").Value = "'531.51"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("E38").Value = "  +2.63%  "
$ws.Range("D39").Value = "'38.25"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").Value = "0.0₃0802"
$ws.Range("E40").Value = "  +4.91%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "3.621.59"
$ws.Range("E42").Value = "  +9.39%  "
$ws.Range("D43").Value = "'0.138"
$ws.Range("E43").Value = "  +4.37%  "
$ws.Range("D44").Value = "'3.43"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("D45").Value = "'0.0466"
$ws.Range("E45").Value = "  +5.03%  "
$ws.Range("D46").Value = "'3.45"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").Value = "'2.91"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "'9.31"
$ws.Range("E48").Value = "  +4.68%  "
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'1.45"
$ws.Range("E51").Value = "  +6.77%  "

# Reset style for cells where an apostrophe (text) prefix was used,
# so no extra quotePrefix style lingers on the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"

Write-Output "Applied cryptos update"
